# Add two new trailing columns (AI, AJ) to the Balance Sheet extract sheet:
#   AI -> other_assets_plug
#   AJ -> other_liabilities_plug
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers, styled the same as the rest of the header row (bold/centered/bordered).
$ws.Range("AI1").Value = "other_assets_plug"
$ws.Range("AJ1").Value = "other_liabilities_plug"
$ws.Range("AG1").Copy()
$ws.Range("AI1:AJ1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AI2").Value = -480373
$ws.Range("AJ2").Value = -588106
$ws.Range("AI3").Value = 3390
$ws.Range("AJ3").Value = -10468
$ws.Range("AI4").Value = 388322
$ws.Range("AJ4").Value = -23279782
$ws.Range("AI5").Value = 1013241
$ws.Range("AJ5").Value = 146092
$ws.Range("AI6").Value = 60257934
$ws.Range("AJ6").Value = 12862123
$ws.Range("AI7").Value = 3837
$ws.Range("AJ7").Value = 11942
$ws.Range("AI8").Value = 4688209
$ws.Range("AJ8").Value = 18752997
$ws.Range("AI9").Value = 1545113804
$ws.Range("AJ9").Value = 783324574
$ws.Range("AI10").Value = -63747719
$ws.Range("AJ10").Value = -116044425
$ws.Range("AI11").Value = 7833024
$ws.Range("AJ11").Value = 5999293
$ws.Range("AI12").Value = 2800303
$ws.Range("AJ12").Value = -514739

# Row 13 (OHIO_STATE_UNIVERSITY_THE) has no data for these two columns, same
# as the other blank cells already present in that row - leave AI13/AJ13 unset.

$ws.Range("AI14").Value = 64890640
$ws.Range("AJ14").Value = 1745473
$ws.Range("AI15").Value = 19794337
$ws.Range("AJ15").Value = 3506904
$ws.Range("AI16").Value = 35006
$ws.Range("AJ16").Value = 2986
$ws.Range("AI17").Value = -1770749
$ws.Range("AJ17").Value = 104868
$ws.Range("AI18").Value = 9970455392
$ws.Range("AJ18").Value = -1811884569

# Row 19 (UNIVERSITY_OF_COLORADO) is likewise blank for these two columns.

$ws.Range("AI20").Value = 387397
$ws.Range("AJ20").Value = -171441
